$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Ror1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1599003333333333
$ws.Range("H2").Value = 0.479701
$ws.Range("I2").Value = 0.0264777194346773
$ws.Range("J2").Value = 0.02647771943467731
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.444123
$ws.Range("N2").Value = 1.332369
$ws.Range("O2").Value = 0.02960111678165545
$ws.Range("P2").Value = 0.02960111678165545
$ws.Range("Q2").Value = 0.071015415741
$ws.Range("R2").Value = 0.639138741669
$ws.Range("S2").Value = 0.0007837700650977911
$ws.Range("T2").Value = 0.000783770065097791

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Ror1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1599003333333333
$ws.Range("H3").Value = 0.479701
$ws.Range("I3").Value = 0.0264777194346773
$ws.Range("J3").Value = 0.02647771943467731
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 12.606804
$ws.Range("N3").Value = 37.820412
$ws.Range("O3").Value = 0.8402525369040582
$ws.Range("P3").Value = 0.8402525369040581
$ws.Range("Q3").Value = 2.015832161868
$ws.Range("R3").Value = 18.142489456812
$ws.Range("S3").Value = 0.02224797092642149
$ws.Range("T3").Value = 0.02224797092642149

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Ror1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1599003333333333
$ws.Range("H4").Value = 0.479701
$ws.Range("I4").Value = 0.0264777194346773
$ws.Range("J4").Value = 0.02647771943467731
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.952662333333333
$ws.Range("N4").Value = 5.857987
$ws.Range("O4").Value = 0.1301463463142864
$ws.Range("P4").Value = 0.1301463463142864
$ws.Range("Q4").Value = 0.3122313579874444
$ws.Range("R4").Value = 2.810082221887
$ws.Range("S4").Value = 0.003445978443158025
$ws.Range("T4").Value = 0.003445978443158024

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Ror1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.879152
$ws.Range("H5").Value = 17.637456
$ws.Range("I5").Value = 0.9735222805653226
$ws.Range("J5").Value = 0.9735222805653228
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.444123
$ws.Range("N5").Value = 1.332369
$ws.Range("O5").Value = 0.02960111678165545
$ws.Range("P5").Value = 0.02960111678165545
$ws.Range("Q5").Value = 2.611066623696
$ws.Range("R5").Value = 23.499599613264
$ws.Range("S5").Value = 0.02881734671655766
$ws.Range("T5").Value = 0.02881734671655766

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Ror1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.879152
$ws.Range("H6").Value = 17.637456
$ws.Range("I6").Value = 0.9735222805653226
$ws.Range("J6").Value = 0.9735222805653228
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.606804
$ws.Range("N6").Value = 37.820412
$ws.Range("O6").Value = 0.8402525369040582
$ws.Range("P6").Value = 0.8402525369040581
$ws.Range("Q6").Value = 74.11731695020799
$ws.Range("R6").Value = 667.0558525518719
$ws.Range("S6").Value = 0.8180045659776367
$ws.Range("T6").Value = 0.8180045659776367

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Ror1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.879152
$ws.Range("H7").Value = 17.637456
$ws.Range("I7").Value = 0.9735222805653226
$ws.Range("J7").Value = 0.9735222805653228
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.952662333333333
$ws.Range("N7").Value = 5.857987
$ws.Range("O7").Value = 0.1301463463142864
$ws.Range("P7").Value = 0.1301463463142864
$ws.Range("Q7").Value = 11.47999866234133
$ws.Range("R7").Value = 103.319987961072
$ws.Range("S7").Value = 0.1267003678711284
$ws.Range("T7").Value = 0.1267003678711284

